$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88").Value = "sqrt(x**2 + y**2) - 1.0"
$ws.Range("B88").Value = 20000
$ws.Range("C88").Value = 1
$ws.Range("D88").Value = 1
$ws.Range("E88").Value = 0.00001008981780614704
$ws.Range("F88").Value = 28.36333250999451
$ws.Range("G88").Value = 705.1357590985656

# H88/I88 must hold the literal text "False" (not a Boolean). Assigning the
# string "False"/"True" straight to .Value auto-coerces to a Boolean (same
# as real Excel), so build it as a text formula result in a scratch cell and
# paste-special the value back in as plain text.
$ws.Range("Z1").Formula = "=""False"""
$ws.Range("Z1").Copy()
$ws.Range("H88").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("I88").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
